$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 3: clear D3, E3, F3
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4: clear E4, F4 (D4 stays 20)
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5: clear D5, F5; set E5 to 4
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 4
$ws.Range("F5").ClearContents()

# Update selection to E5
$ws.Activate()
$ws.Range("E5").Select()
